$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: URL in B17 (hyperlink style), description in C17
$ws.Range("B17").Value = "https://git-scm.com/cheat-sheet"
$ws.Range("C17").Value = "Comados Básicos de Git"

# Add hyperlink for B17
$ws.Hyperlinks.Add($ws.Range("B17"), "https://git-scm.com/cheat-sheet", "", "", "https://git-scm.com/cheat-sheet")

# Copy style from B16/C16 to B17/C17 to keep consistent formatting
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)

# Update selection
$ws.Range("C18").Select()

$wb.Save()
